# "Generate Report for Handback" — refresh the handoff/handback timestamps
# for the 7b57cec4-62d6-4afc-b527-bbd37cd5c6e1.md file across the Overview,
# zh-cn and de-de sheets of the handback-status report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview!G3 — "Latest HO Xliff Generate Date" for the de-de column of the
# 7b57cec4-62d6-4afc-b527-bbd37cd5c6e1.md row.
$overview.Range("G3").Value = "2016-08-12 23:01:41"

# zh-cn sheet, row 3 (7b57cec4-62d6-4afc-b527-bbd37cd5c6e1 file):
#   H = Correspond Handoff Datetime
#   K = Correspond Handback DateTime
$zhcn.Range("H3").Value = "2016-08-12 23:01:34"
$zhcn.Range("K3").Value = "2016-08-12 23:02:18"

# de-de sheet, row 3 (7b57cec4-62d6-4afc-b527-bbd37cd5c6e1 file):
#   H = Correspond Handoff Datetime
#   K = Correspond Handback DateTime
$dede.Range("H3").Value = "2016-08-12 23:01:41"
$dede.Range("K3").Value = "2016-08-12 23:02:27"
